# Apply updated crypto prices/volumes/rankings to the worksheet (cryptos.xlsx)
# Mirrors the per-row data refresh performed by the GitHub Actions scraper.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '62.922.32'
$ws.Range('E2').Value = '  -5.26%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '3.114.74'
$ws.Range('E3').Value = '  -5.73%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.08%  '

# Row 5: BNB
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '561.91'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.24%  '

# Row 6: Solana
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -9.88%  '

# Row 7: USDC
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.06%  '

# Row 8: XRP
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -9.95%  '

# Row 9: LidoStakedEther
$ws.Range('D9').Value = '3.112.21'
$ws.Range('E9').Value = '  -5.80%  '

# Row 10: Toncoin
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.70'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.17%  '

# Row 11: Dogecoin
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.116'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.41%  '

# Row 12: Cardano
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.379'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.98%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range('D13').Value = '3.667.81'
$ws.Range('E13').Value = '  -5.66%  '

# Row 14: TRON
$ws.Range('E14').Value = '  -2.12%  '

# Row 15: WrappedBTC
$ws.Range('D15').Value = '63.054.11'
$ws.Range('E15').Value = '  -5.09%  '

# Row 16: Avalanche
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.71%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '3.130.75'
$ws.Range('E17').Value = '  -5.79%  '

# Row 18: ShibaInu
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000154'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.06%  '

# Row 19: BitcoinCash
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '406.47'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.95%  '

# Row 20: Chainlink
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.26%  '

# Row 21: Polkadot
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.82%  '

# Row 22: Uniswap
$ws.Range('E22').Value = '  -3.91%  '

# Row 23: Dai
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.18%  '

# Row 24: LEO
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.02%  '

# Row 25: Litecoin
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.98%  '

# Row 26: Kaspa
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.201'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.95%  '

# Row 27: Polygon
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.489'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.90%  '

# Row 28: PEPE
$ws.Range('E28').Value = '  -11.48%  '

# Row 29: InternetComputer(DFINITY)
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.70'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.46%  '

# Row 30: Binance-PegBSC-USD
$ws.Range('E30').Value = '  +0.28%  '

# Row 31: USDe
$ws.Range('E31').Value = '  -0.08%  '

# Row 32: PancakeSwap
$ws.Range('E32').Value = '  -6.67%  '

# Row 33: EthereumClassic
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.83%  '

# Row 34: NEARProtocol
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.87'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.72%  '

# Row 35: Aptos
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.24'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.36%  '

# Row 36: Monero
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '154.84'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.23%  '

# Row 37: Fetch.AI
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.10'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.33%  '

# Row 38: ImmutableX
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.33'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.72%  '

# Row 39: Maker
$ws.Range('D39').Value = '2.693.27'
$ws.Range('E39').Value = '  -5.68%  '

# Row 40: Stacks
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.65'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -8.31%  '

# Row 41: EnergySwap
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '23.65'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.98%  '

# Row 42: Filecoin
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.11'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.99%  '

# Row 43: OKB
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.02%  '

# Row 44: Mantle
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.694'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.81%  '

# Row 45: Hedera
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0607'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.63%  '

# Row 46: VeChain
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0256'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.88%  '

# Row 47: RenderToken
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.19'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -11.95%  '

# Row 48: InjectiveProtocol
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.88'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -9.43%  '

# Row 49: Bittensor
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '281.03'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -9.50%  '

# Row 50: FirstDigitalUSD
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.999'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.09%  '

# Row 51: Stellar
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0971'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -6.48%  '
